$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New pruebab2b.com URLs replacing the old plastico.com ones.
$ws.Range("B32").Value = "http://www.pruebab2b.com/contactenos"
$ws.Range("B33").Value = "http://www.pruebab2b.com/quienes-somos"
$ws.Range("B34").Value = "http://www.pruebab2b.com/aviso-de-privacidad"
$ws.Range("B35").Value = "http://www.pruebab2b.com/condiciones-del-servicio"
$ws.Range("B36").Value = "http://www.pruebab2b.com/copyright"

# B32:B34 pick up the hyperlink-style formatting that B35/B36 already had.
$ws.Range("B32").Style = "Hipervínculo"
$ws.Range("B33").Style = "Hipervínculo"
$ws.Range("B34").Style = "Hipervínculo"

# B32:B34 had no hyperlink objects before; add them now.
$ws.Hyperlinks.Add($ws.Range("B32"), "http://www.pruebab2b.com/contactenos")
$ws.Hyperlinks.Add($ws.Range("B33"), "http://www.pruebab2b.com/quienes-somos")
$ws.Hyperlinks.Add($ws.Range("B34"), "http://www.pruebab2b.com/aviso-de-privacidad")

# B35/B36 already had hyperlink objects pointing at the old plastico.com
# addresses; repoint them at the new pruebab2b.com addresses.
$ws.Hyperlinks.Item(1).Address = "http://www.pruebab2b.com/condiciones-del-servicio"
$ws.Hyperlinks.Item(2).Address = "http://www.pruebab2b.com/copyright"

# Selection moved from C34 to B32.
$ws.Range("B32").Select()
